$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'286.33"
$ws.Range("E2").Value = "'0.64%"
$ws.Range("G2").Value = "'2"

# Row 3
$ws.Range("D3").Value = "'29.34"
$ws.Range("E3").Value = "'3.60%"
$ws.Range("G3").Value = "'2"

# Row 4
$ws.Range("D4").Value = "'5.079"
$ws.Range("E4").Value = "'3.54%"
$ws.Range("G4").Value = "'2"

# Row 5
$ws.Range("D5").Value = "'0.06657"
$ws.Range("E5").Value = "'1.29%"
$ws.Range("G5").Value = "'2"

# Row 6
$ws.Range("D6").Value = "'7.327"
$ws.Range("E6").Value = "'1.61%"
$ws.Range("G6").Value = "'2"

# Row 7
$ws.Range("D7").Value = "'3.405"
$ws.Range("E7").Value = "'0.88%"
$ws.Range("G7").Value = "'2"

# Row 8
$ws.Range("D8").Value = "'1.352"
$ws.Range("E8").Value = "'-2.21%"
$ws.Range("G8").Value = "'2"

# Row 9
$ws.Range("D9").Value = "'0.9235"
$ws.Range("E9").Value = "'1.05%"
$ws.Range("G9").Value = "'2"

# Row 10
$ws.Range("D10").Value = "'0.1565"
$ws.Range("E10").Value = "'0.39%"
$ws.Range("G10").Value = "'2"

# Row 11
$ws.Range("D11").Value = "'0.06326"
$ws.Range("E11").Value = "'0.79%"
$ws.Range("G11").Value = "'2"

# Row 12
$ws.Range("D12").Value = "'0.07561"
$ws.Range("E12").Value = "'-0.54%"
$ws.Range("G12").Value = "'2"

# Row 13
$ws.Range("D13").Value = "'0.02897"
$ws.Range("E13").Value = "'-0.94%"
$ws.Range("G13").Value = "'2"

# Row 14
$ws.Range("D14").Value = "'0.08993"
$ws.Range("E14").Value = "'0.36%"
$ws.Range("G14").Value = "'2"

# Row 15
$ws.Range("D15").Value = "'0.001594"
$ws.Range("E15").Value = "'-0.23%"
$ws.Range("G15").Value = "'2"

# Row 16
$ws.Range("E16").Value = "'1.01%"
$ws.Range("G16").Value = "'2"

# Row 17
$ws.Range("D17").Value = "'0.0006483"
$ws.Range("E17").Value = "'-0.39%"
$ws.Range("G17").Value = "'2"

# Row 18
$ws.Range("D18").Value = "'0.006256"
$ws.Range("E18").Value = "'4.55%"
$ws.Range("G18").Value = "'2"

# Row 19
$ws.Range("D19").Value = "'3.457"
$ws.Range("E19").Value = "'-1.00%"
$ws.Range("G19").Value = "'2"

# Row 20
$ws.Range("D20").Value = "'2.229"
$ws.Range("E20").Value = "'-0.42%"
$ws.Range("G20").Value = "'2"

# Row 21
$ws.Range("D21").Value = "'0.3212"
$ws.Range("E21").Value = "'0.88%"
$ws.Range("G21").Value = "'2"

# Row 22
$ws.Range("D22").Value = "'0.1309"
$ws.Range("E22").Value = "'-2.72%"
$ws.Range("G22").Value = "'2"

# Row 23
$ws.Range("D23").Value = "'4.069"
$ws.Range("E23").Value = "'2.81%"
$ws.Range("G23").Value = "'2"

# Row 24
$ws.Range("D24").Value = "'0.1550"
$ws.Range("E24").Value = "'2.24%"
$ws.Range("G24").Value = "'2"

# Row 25
$ws.Range("D25").Value = "'0.001193"
$ws.Range("E25").Value = "'0.84%"
$ws.Range("G25").Value = "'2"

# Row 26
$ws.Range("E26").Value = "'-5.56%"
$ws.Range("G26").Value = "'2"

# Row 27
$ws.Range("D27").Value = "'0.0001249"
$ws.Range("E27").Value = "'6.15%"
$ws.Range("G27").Value = "'2"

# Row 28
$ws.Range("D28").Value = "'0.0001617"
$ws.Range("E28").Value = "'-1.67%"
$ws.Range("G28").Value = "'2"

# Row 29
$ws.Range("G29").Value = "'2"

# Row 30
$ws.Range("G30").Value = "'2"

# Row 31
$ws.Range("G31").Value = "'2"

# Row 32
$ws.Range("G32").Value = "'2"

# Row 33
$ws.Range("G33").Value = "'2"

# Row 34
$ws.Range("G34").Value = "'2"

# Row 35
$ws.Range("G35").Value = "'2"

# Row 36
$ws.Range("G36").Value = "'2"

# Row 37
$ws.Range("G37").Value = "'2"

# Row 38
$ws.Range("G38").Value = "'2"

# Row 39
$ws.Range("G39").Value = "'2"

# Row 40
$ws.Range("E40").Value = "'0.71%"
$ws.Range("G40").Value = "'2"

# Row 41
$ws.Range("D41").Value = "'0.006712"
$ws.Range("E41").Value = "'-1.74%"
$ws.Range("G41").Value = "'2"

# Row 42
$ws.Range("D42").Value = "'0.1238"
$ws.Range("E42").Value = "'-12.39%"
$ws.Range("G42").Value = "'2"

# Row 43
$ws.Range("D43").Value = "'0.001979"
$ws.Range("E43").Value = "'-2.74%"
$ws.Range("G43").Value = "'2"

# Row 44
$ws.Range("D44").Value = "'0.01251"
$ws.Range("E44").Value = "'6.83%"
$ws.Range("G44").Value = "'2"

# Row 45
$ws.Range("D45").Value = "'0.00005601"
$ws.Range("E45").Value = "'0.97%"
$ws.Range("G45").Value = "'2"

# Row 46
$ws.Range("D46").Value = "'0.01306"
$ws.Range("E46").Value = "'-29.17%"
$ws.Range("G46").Value = "'2"

# Row 47
$ws.Range("E47").Value = "'20.85%"
$ws.Range("G47").Value = "'2"

# Row 48
$ws.Range("G48").Value = "'2"

# Row 49
$ws.Range("G49").Value = "'2"

# Row 50
$ws.Range("G50").Value = "'2"

# Row 51
$ws.Range("G51").Value = "'2"

Write-Output "applied 116 cell updates"
